$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 165
    $ws.Range("F3").Value = 33
    $ws.Range("F4").Value = 116
}
